$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "+211.03 ± 0.29"
$ws.Range("D2").Value = "+0.07 ± 0.00"
$ws.Range("E2").Value = "+204.50 ± 0.29"
$ws.Range("F2").Value = "+6.26 ± 0.02"
$ws.Range("G2").Value = "+0.56 ± 0.01"
$ws.Range("B4").Value = "+4.57 ± 0.02"
$ws.Range("E4").Value = "+3.21 ± 0.01"
$ws.Range("F4").Value = "+0.97 ± 0.01"
$ws.Range("G4").Value = "+0.32 ± 0.01"
$ws.Range("B5").Value = "+5.24 ± 0.01"
$ws.Range("F5").Value = "+5.09 ± 0.01"
$ws.Range("G5").Value = "+0.15 ± 0.00"
$ws.Range("E6").Value = "+0.07 ± 0.00"
$ws.Range("B8").Value = "+201.22 ± 0.29"
$ws.Range("E8").Value = "+201.22 ± 0.29"
$ws.Range("B9").Value = "-211.13 ± 0.29"
$ws.Range("D9").Value = "-0.07 ± 0.00"
$ws.Range("E9").Value = "-204.57 ± 0.29"
$ws.Range("F9").Value = "-6.29 ± 0.02"
$ws.Range("G9").Value = "-0.56 ± 0.01"
$ws.Range("D10").Value = "-0.07 ± 0.00"
$ws.Range("B12").Value = "-149.54 ± 0.27"
$ws.Range("E12").Value = "-146.04 ± 0.27"
$ws.Range("F12").Value = "-3.26 ± 0.01"
$ws.Range("G12").Value = "-0.24 ± 0.00"
$ws.Range("B13").Value = "-10.57 ± 0.02"
$ws.Range("E13").Value = "-10.32 ± 0.02"
$ws.Range("F13").Value = "-0.24 ± 0.00"
$ws.Range("B14").Value = "-48.98 ± 0.02"
$ws.Range("E14").Value = "-46.04 ± 0.02"
$ws.Range("F14").Value = "-2.67 ± 0.01"
$ws.Range("G14").Value = "-0.28 ± 0.00"
$ws.Range("E15").Value = "-1.89 ± 0.00"
$ws.Range("D17").Value = "+0.09 ± 0.00"
$ws.Range("F17").Value = "+6.28 ± 0.05"
$ws.Range("G17").Value = "+5.78 ± 0.15"
